$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To-Do-List")

# D16: Week value for "Implement Security" row
$ws.Range("D16").Value = 7

# Row 29: new task "Chapters page" entry
$ws.Range("B29").Value = "Chapters page"
$ws.Range("D29").Value = 7
$ws.Range("E29").Value = "Pending"

# Row 30: new task "task list"
$ws.Range("B30").Value = "task list"
$ws.Range("C30").Value = "added to my blue script page"
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = "Not Started"

# Row 31: new task "create slideshow"
$ws.Range("B31").Value = "create slideshow"
$ws.Range("C31").Value = "for home page"
$ws.Range("D31").Value = 7
$ws.Range("E31").Value = "Pending"

# Remove the stats note at C33
$ws.Range("C33").Clear()

# Update selection / view
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("E31").Select() | Out-Null
